# Update the Trade_Data worksheet:
#  - shrink the trade list from 120 rows (119 trades) down to 89 rows (88 trades)
#  - refresh the "symbol" column (C) for every remaining trade row
#  - flip "signal" (B) from BUY to SELL for the tail of the list (rows 67-89)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trade_Data")

# 1) Drop the now-unused trailing rows (90:120) so the sheet's used range
#    becomes A1:C89, matching the trimmed trade list.
$ws.Range("A90:C120").EntireRow.Delete()

# 2) New symbol values for every data row (rows 2-89, i.e. trades 0-87).
$symbols = @(
    "ADX","AMH","APO","BC","BDN","BHF","CDNS","CG","CLW","CTXS",
    "DHI","DOOR","EGP","ELAN","EWBC","EXC","EXLS","FE","FR","GRBK",
    "HE","HR","IBP","INTC","ITGR","IVR","JBGS","JKHY","KBH","KO",
    "LEN","LMT","LNT","LQDA","MDU","MRCC","MU","MYE","NDSN","NOVA",
    "NX","OPI","ORC","PDM","PEG","PEP","PGRE","PHR","PLXS","PRIM",
    "PTC","RMBS","SLG","SO","SRE","SSD","STWD","TWO","UMH","UNP",
    "VNO","WDC","WRB","XAN","YMAB","AAWW","ALLO","BCEI","BPMC","BTU",
    "CEIX","COG","CRC","DO","ENTA","EQT","GTHX","LBRT","PEB","PUMP",
    "PVAC","RLI","SCHN","SGMO","UNFI","VIRT","VYGR","X"
)

$symbolRange = New-Object 'object[,]' $symbols.Length,1
for ($i = 0; $i -lt $symbols.Length; $i++) {
    $symbolRange[$i,0] = $symbols[$i]
}
$ws.Range("C2:C89").Value = $symbolRange

# 3) Rows 67-89 (trades 65-87) switch their signal from BUY to SELL.
$signals = New-Object 'object[,]' 23,1
for ($i = 0; $i -lt 23; $i++) {
    $signals[$i,0] = "SELL"
}
$ws.Range("B67:B89").Value = $signals

Write-Output "Trade_Data updated: dimension trimmed to A1:C89, symbols refreshed, tail signals set to SELL."
